# Auto-generated Excel COM-interop script
# Applies scheduled-runner value updates to Belias_Profits workbook sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 5537.839
$ws.Range("I106").Value = 4062.1428
$ws.Range("K106").Value = 4062.1428
$ws.Range("M106").Value = -3431.1428
$ws.Range("H137").Value = 618726.4399999999
$ws.Range("I137").Value = 1695.5483
$ws.Range("J137").Value = 942929.1
$ws.Range("K137").Value = 5086.644899999999
$ws.Range("L137").Value = 2828787.3
$ws.Range("M137").Value = -2536.644899999999
$ws.Range("N137").Value = -2833887.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2966.397
$ws.Range("I32").Value = 2479.2742
$ws.Range("J32").Value = 8000
$ws.Range("K32").Value = 2479.2742
$ws.Range("L32").Value = 8000
$ws.Range("M32").Value = -2192.2742
$ws.Range("N32").Value = -8574
$ws.Range("H40").Value = 9920
$ws.Range("J40").Value = 9920
$ws.Range("L40").Value = 9920
$ws.Range("N40").Value = -10272
$ws.Range("H132").Value = 2149.394
$ws.Range("I132").Value = 1683.48
$ws.Range("J132").Value = 3605.375
$ws.Range("K132").Value = 5050.440000000001
$ws.Range("L132").Value = 10816.125
$ws.Range("M132").Value = -2520.440000000001
$ws.Range("N132").Value = -15876.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1830.6522
$ws.Range("I86").Value = 1637.1052
$ws.Range("K86").Value = 1637.1052
$ws.Range("M86").Value = -514.1052
$ws.Range("H89").Value = 1830.6522
$ws.Range("I89").Value = 1637.1052
$ws.Range("K89").Value = 8185.526
$ws.Range("M89").Value = -2569.526
$ws.Range("H107").Value = 2975.4
$ws.Range("I107").Value = 3083.7778
$ws.Range("J107").Value = 2000
$ws.Range("K107").Value = 3083.7778
$ws.Range("L107").Value = 2000
$ws.Range("M107").Value = -1163.7778
$ws.Range("N107").Value = -5840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5883379
$ws.Range("I31").Value = 616.1667
$ws.Range("J31").Value = 10205409
$ws.Range("K31").Value = 616.1667
$ws.Range("L31").Value = 10205409
$ws.Range("M31").Value = -321.1667
$ws.Range("N31").Value = -10205999
$ws.Range("H34").Value = 5883379
$ws.Range("I34").Value = 616.1667
$ws.Range("J34").Value = 10205409
$ws.Range("K34").Value = 616.1667
$ws.Range("L34").Value = 10205409
$ws.Range("M34").Value = -414.1667
$ws.Range("N34").Value = -10205813
$ws.Range("H38").Value = 22000
$ws.Range("I38").Value = 10000
$ws.Range("J38").Value = 25000
$ws.Range("K38").Value = 10000
$ws.Range("L38").Value = 25000
$ws.Range("M38").Value = -9623
$ws.Range("N38").Value = -25754
$ws.Range("H42").Value = 5707.5
$ws.Range("I42").Value = 4808.5713
$ws.Range("J42").Value = 12000
$ws.Range("K42").Value = 4808.5713
$ws.Range("L42").Value = 12000
$ws.Range("M42").Value = -4215.5713
$ws.Range("N42").Value = -13186
$ws.Range("H46").Value = 22000
$ws.Range("I46").Value = 10000
$ws.Range("J46").Value = 25000
$ws.Range("K46").Value = 10000
$ws.Range("L46").Value = 25000
$ws.Range("M46").Value = -9789
$ws.Range("N46").Value = -25422
$ws.Range("H86").Value = 10017.917
$ws.Range("I86").Value = 4338.375
$ws.Range("K86").Value = 4338.375
$ws.Range("M86").Value = -3215.375
$ws.Range("H89").Value = 10017.917
$ws.Range("I89").Value = 4338.375
$ws.Range("K89").Value = 21691.875
$ws.Range("M89").Value = -16075.875
$ws.Range("H107").Value = 1540.8422
$ws.Range("I107").Value = 681.5
$ws.Range("J107").Value = 2165.818
$ws.Range("K107").Value = 681.5
$ws.Range("L107").Value = 2165.818
$ws.Range("M107").Value = 1238.5
$ws.Range("N107").Value = -6005.818

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 112256.445
$ws.Range("I2").Value = 23.75
$ws.Range("J2").Value = 202042.6
$ws.Range("K2").Value = 142.5
$ws.Range("L2").Value = 1212255.6
$ws.Range("M2").Value = -29.5
$ws.Range("N2").Value = -1212481.6
$ws.Range("H15").Value = 1198.8
$ws.Range("I15").Value = 45
$ws.Range("J15").Value = 1376.3077
$ws.Range("K15").Value = 135
$ws.Range("L15").Value = 4128.9231
$ws.Range("M15").Value = 5
$ws.Range("N15").Value = -4408.9231
$ws.Range("H39").Value = 3739.8572
$ws.Range("I39").Value = 800
$ws.Range("J39").Value = 4229.8335
$ws.Range("K39").Value = 2400
$ws.Range("L39").Value = 12689.5005
$ws.Range("M39").Value = -2106
$ws.Range("N39").Value = -13277.5005
$ws.Range("H58").Value = 3764.5557
$ws.Range("I58").Value = 1002.5
$ws.Range("J58").Value = 4109.8125
$ws.Range("K58").Value = 3007.5
$ws.Range("L58").Value = 12329.4375
$ws.Range("M58").Value = -2879.5
$ws.Range("N58").Value = -12585.4375
$ws.Range("H64").Value = 4201.857
$ws.Range("I64").Value = 2565.75
$ws.Range("J64").Value = 6383.3335
$ws.Range("K64").Value = 7697.25
$ws.Range("L64").Value = 19150.0005
$ws.Range("M64").Value = -7427.25
$ws.Range("N64").Value = -19690.0005
$ws.Range("H67").Value = 4201.857
$ws.Range("I67").Value = 2565.75
$ws.Range("J67").Value = 6383.3335
$ws.Range("K67").Value = 7697.25
$ws.Range("L67").Value = 19150.0005
$ws.Range("M67").Value = -6761.25
$ws.Range("N67").Value = -21022.0005
$ws.Range("H76").Value = 4920
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").Value = ""
$ws.Range("H79").Value = 4920
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").Value = ""
$ws.Range("H112").Value = 100004070
$ws.Range("I112").Value = 1450
$ws.Range("J112").Value = 125004740
$ws.Range("K112").Value = 4350
$ws.Range("L112").Value = 375014220
$ws.Range("M112").Value = -3242
$ws.Range("N112").Value = -375016436
$ws.Range("H131").Value = 850.7
$ws.Range("I131").Value = 442
$ws.Range("J131").Value = 901.2135
$ws.Range("K131").Value = 1326
$ws.Range("L131").Value = 2703.6405
$ws.Range("M131").Value = 3714
$ws.Range("N131").Value = -12783.6405

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2166.6667
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6840
$ws.Range("H132").Value = 2063.4
$ws.Range("I132").Value = 1852.9032
$ws.Range("J132").Value = 2529.5
$ws.Range("K132").Value = 5558.7096
$ws.Range("L132").Value = 7588.5
$ws.Range("M132").Value = -3028.7096
$ws.Range("N132").Value = -12648.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 49900
$ws.Range("J133").Value = 49900
$ws.Range("L133").Value = 49900
$ws.Range("N133").Value = -54960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 20000
$ws.Range("J42").Value = 20000
$ws.Range("L42").Value = 20000
$ws.Range("N42").Value = -20756
$ws.Range("H107").Value = 17509.182
$ws.Range("I107").Value = 23963.875
$ws.Range("K107").Value = 71891.625
$ws.Range("M107").Value = -69971.625
$ws.Range("H136").Value = 3345.8064
$ws.Range("I136").Value = 3569.5789
$ws.Range("K136").Value = 10708.7367
$ws.Range("M136").Value = -8158.736699999999
